# Updated symbol list on Mon Jan 30 20:33:00 UTC 2023 with GitHub Actions
# Refresh price/volume figures (and fix the BOLO / CoinbaseStockToken row order)
# in the cryptos sheet to match the latest scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'305.64"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'-3.89%"
$ws.Range("E2").Style = "Normal"
# Row 3
$ws.Range("D3").Value = "'37.02"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'-7.03%"
$ws.Range("E3").Style = "Normal"
# Row 4
$ws.Range("D4").Value = "'5.104"
$ws.Range("D4").Style = "Normal"
# Row 5
$ws.Range("D5").Value = "'0.07705"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'-6.32%"
$ws.Range("E5").Style = "Normal"
# Row 6
$ws.Range("D6").Value = "'4.379"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'1.08%"
$ws.Range("E6").Style = "Normal"
# Row 7
$ws.Range("D7").Value = "'8.209"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'-1.72%"
$ws.Range("E7").Style = "Normal"
# Row 8
$ws.Range("D8").Value = "'1.869"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'-9.70%"
$ws.Range("E8").Style = "Normal"
# Row 9
$ws.Range("E9").Value = "'-4.67%"
$ws.Range("E9").Style = "Normal"
# Row 10
$ws.Range("D10").Value = "'0.9189"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'-2.14%"
$ws.Range("E10").Style = "Normal"
# Row 11
$ws.Range("D11").Value = "'0.1223"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'-10.81%"
$ws.Range("E11").Style = "Normal"
# Row 12
$ws.Range("D12").Value = "'0.1896"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'-4.02%"
$ws.Range("E12").Style = "Normal"
# Row 13
$ws.Range("D13").Value = "'0.08741"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'-4.05%"
$ws.Range("E13").Style = "Normal"
# Row 14
$ws.Range("D14").Value = "'0.03404"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'-3.69%"
$ws.Range("E14").Style = "Normal"
# Row 15
$ws.Range("D15").Value = "'0.09691"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'-1.14%"
$ws.Range("E15").Style = "Normal"
# Row 16
$ws.Range("D16").Value = "'0.001366"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'-3.59%"
$ws.Range("E16").Style = "Normal"
# Row 17
$ws.Range("D17").Value = "'0.006094"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'-2.69%"
$ws.Range("E17").Style = "Normal"
# Row 18
$ws.Range("E18").Value = "'-3.53%"
$ws.Range("E18").Style = "Normal"
# Row 19
$ws.Range("E19").Value = "'-3.56%"
$ws.Range("E19").Style = "Normal"
# Row 20
$ws.Range("D20").Value = "'0.1284"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'-2.33%"
$ws.Range("E20").Style = "Normal"
# Row 21
$ws.Range("D21").Value = "'5.030"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'1.59%"
$ws.Range("E21").Style = "Normal"
# Row 22
$ws.Range("D22").Value = "'0.2503"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'2.22%"
$ws.Range("E22").Style = "Normal"
# Row 23
$ws.Range("D23").Value = "'0.02116"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'5,191.14%"
$ws.Range("E23").Style = "Normal"
# Row 24
$ws.Range("D24").Value = "'0.04329"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'-0.51%"
$ws.Range("E24").Style = "Normal"
# Row 25
$ws.Range("D25").Value = "'0.001219"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'-0.85%"
$ws.Range("E25").Style = "Normal"
# Row 26
$ws.Range("D26").Value = "'0.004464"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'-7.49%"
$ws.Range("E26").Style = "Normal"
# Row 27
$ws.Range("D27").Value = "'0.0001356"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'4.45%"
$ws.Range("E27").Style = "Normal"
# Row 39
$ws.Range("D39").Value = "'0.02225"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'0.27%"
$ws.Range("E39").Style = "Normal"
# Row 40
$ws.Range("D40").Value = "'0.04898"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'-6.22%"
$ws.Range("E40").Style = "Normal"
# Row 41
$ws.Range("D41").Value = "'0.007611"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'-1.87%"
$ws.Range("E41").Style = "Normal"
# Row 42
$ws.Range("D42").Value = "'0.009938"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'2.67%"
$ws.Range("E42").Style = "Normal"
# Row 43
$ws.Range("E43").Value = "'-5.27%"
$ws.Range("E43").Style = "Normal"
# Row 44
$ws.Range("D44").Value = "'0.002003"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'-2.10%"
$ws.Range("E44").Style = "Normal"
# Row 45
$ws.Range("D45").Value = "'0.008817"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'-8.72%"
$ws.Range("E45").Style = "Normal"
# Row 46
$ws.Range("D46").Value = "'0.00006966"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'5.40%"
$ws.Range("E46").Style = "Normal"
# Row 47
$ws.Range("E47").Value = "'0.72%"
$ws.Range("E47").Style = "Normal"
# Row 48
$ws.Range("B48").Value = "CoinbaseStockToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/_ZA6fIr53+coinbasestocktoken-coin"
$ws.Range("D48").Value = "'0.001308"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'-22.60%"
$ws.Range("E48").Style = "Normal"
# Row 49
$ws.Range("B49").Value = "BOLO"
$ws.Range("C49").Value = "https://coinranking.com/coin/ogrGe0dEab+bolo-bolo"
$ws.Range("D49").Value = "'0.003009"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'2.37%"
$ws.Range("E49").Style = "Normal"
# Row 50
$ws.Range("D50").Value = "'0.00002113"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'0.72%"
$ws.Range("E50").Style = "Normal"
# Row 51
$ws.Range("D51").Value = "'0.0002012"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'0.72%"
$ws.Range("E51").Style = "Normal"
